# Few test cases added for att.com
#
# Sheet "boaSignup": append 5 new rows of OnlineId/Password style test
# data in columns C:D starting at row 11 (two pre-existing pairs plus
# three brand-new "...2" variants).
#
# Sheet "loginNegativeTest": drop the now-unused "ErrContains" column C,
# and drop the two extra negative-test rows (3 and 4), leaving just the
# header + one data row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("boaSignup")
$ws1.Range("C11").Value = "FGHIJK"
$ws1.Range("D11").Value = "FHG1234$"
$ws1.Range("C12").Value = "LMNOPQ"
$ws1.Range("D12").Value = "pqrst123#"
$ws1.Range("C13").Value = "ABCDE2"
$ws1.Range("D13").Value = "abcd123#"
$ws1.Range("C14").Value = "FGHIJK2"
$ws1.Range("D14").Value = "FHG1234$"
$ws1.Range("C15").Value = "LMNO2"
$ws1.Range("D15").Value = "pqrst123#"

$ws2 = $wb.Worksheets.Item("loginNegativeTest")
$ws2.Range("C1:C4").EntireColumn.Delete()
$ws2.Range("A3:A4").EntireRow.Delete()

$ws1.Activate()
$ws1.Range("C11:D15").Select()

$ws2.Activate()
$ws2.Range("B1").Select()
